# Auto-generated: update leve-profit market-data cells per scheduled-runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 7018.778  # H48
$ws.Cells.Item(48, 10).Value = 7018.778  # J48
$ws.Cells.Item(48, 12).Value = 21056.334  # L48
$ws.Cells.Item(48, 14).Value = -21640.334  # N48
$ws.Cells.Item(56, 8).Value = 7018.778  # H56
$ws.Cells.Item(56, 10).Value = 7018.778  # J56
$ws.Cells.Item(56, 12).Value = 21056.334  # L56
$ws.Cells.Item(56, 14).Value = -22124.334  # N56
$ws.Cells.Item(74, 8).Value = 3997955.5  # H74
$ws.Cells.Item(74, 9).Value = 4724038.5  # I74
$ws.Cells.Item(74, 10).Value = 4500  # J74
$ws.Cells.Item(74, 11).Value = 4724038.5  # K74
$ws.Cells.Item(74, 12).Value = 4500  # L74
$ws.Cells.Item(74, 13).Value = -4723102.5  # M74
$ws.Cells.Item(74, 14).Value = -6372  # N74
$ws.Cells.Item(77, 8).Value = 3997955.5  # H77
$ws.Cells.Item(77, 9).Value = 4724038.5  # I77
$ws.Cells.Item(77, 10).Value = 4500  # J77
$ws.Cells.Item(77, 11).Value = 23620192.5  # K77
$ws.Cells.Item(77, 12).Value = 22500  # L77
$ws.Cells.Item(77, 13).Value = -23615512.5  # M77
$ws.Cells.Item(77, 14).Value = -31860  # N77
$ws.Cells.Item(86, 8).Value = 43481616  # H86
$ws.Cells.Item(86, 9).Value = 3601  # I86
$ws.Cells.Item(86, 10).Value = 58826796  # J86
$ws.Cells.Item(86, 11).Value = 3601  # K86
$ws.Cells.Item(86, 12).Value = 58826796  # L86
$ws.Cells.Item(86, 13).Value = -2478  # M86
$ws.Cells.Item(86, 14).Value = -58829042  # N86
$ws.Cells.Item(89, 8).Value = 43481616  # H89
$ws.Cells.Item(89, 9).Value = 3601  # I89
$ws.Cells.Item(89, 10).Value = 58826796  # J89
$ws.Cells.Item(89, 11).Value = 18005  # K89
$ws.Cells.Item(89, 12).Value = 294133980  # L89
$ws.Cells.Item(89, 13).Value = -12389  # M89
$ws.Cells.Item(89, 14).Value = -294145212  # N89
$ws.Cells.Item(135, 8).Value = 1165.0571  # H135
$ws.Cells.Item(135, 9).Value = 666.1070999999999  # I135
$ws.Cells.Item(135, 10).Value = 3160.8572  # J135
$ws.Cells.Item(135, 11).Value = 5994.9639  # K135
$ws.Cells.Item(135, 12).Value = 28447.7148  # L135
$ws.Cells.Item(135, 13).Value = -3459.9639  # M135
$ws.Cells.Item(135, 14).Value = -33517.7148  # N135
$ws.Cells.Item(138, 8).Value = 1575.49  # H138
$ws.Cells.Item(138, 9).Value = 619.67145  # I138
$ws.Cells.Item(138, 10).Value = 3805.7334  # J138
$ws.Cells.Item(138, 11).Value = 1859.01435  # K138
$ws.Cells.Item(138, 12).Value = 11417.2002  # L138
$ws.Cells.Item(138, 13).Value = 3280.98565  # M138
$ws.Cells.Item(138, 14).Value = -21697.2002  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1485.0869  # H61
$ws.Cells.Item(61, 9).Value = 956.2  # I61
$ws.Cells.Item(61, 10).Value = 2476.75  # J61
$ws.Cells.Item(61, 11).Value = 956.2  # K61
$ws.Cells.Item(61, 12).Value = 2476.75  # L61
$ws.Cells.Item(61, 13).Value = -744.2  # M61
$ws.Cells.Item(61, 14).Value = -2900.75  # N61
$ws.Cells.Item(74, 8).Value = 1020.95  # H74
$ws.Cells.Item(74, 9).Value = 968.069  # I74
$ws.Cells.Item(74, 10).Value = 1160.3636  # J74
$ws.Cells.Item(74, 11).Value = 968.069  # K74
$ws.Cells.Item(74, 12).Value = 1160.3636  # L74
$ws.Cells.Item(74, 13).Value = -94.06899999999996  # M74
$ws.Cells.Item(74, 14).Value = -2908.3636  # N74
$ws.Cells.Item(77, 8).Value = 1020.95  # H77
$ws.Cells.Item(77, 9).Value = 968.069  # I77
$ws.Cells.Item(77, 10).Value = 1160.3636  # J77
$ws.Cells.Item(77, 11).Value = 4840.344999999999  # K77
$ws.Cells.Item(77, 12).Value = 5801.817999999999  # L77
$ws.Cells.Item(77, 13).Value = -472.3449999999993  # M77
$ws.Cells.Item(77, 14).Value = -14537.818  # N77
$ws.Cells.Item(136, 8).Value = 1485.0869  # H136
$ws.Cells.Item(136, 9).Value = 956.2  # I136
$ws.Cells.Item(136, 10).Value = 2476.75  # J136
$ws.Cells.Item(136, 11).Value = 2868.6  # K136
$ws.Cells.Item(136, 12).Value = 7430.25  # L136
$ws.Cells.Item(136, 13).Value = -318.6000000000004  # M136
$ws.Cells.Item(136, 14).Value = -12530.25  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(59, 8).Value = 37500  # H59
$ws.Cells.Item(59, 10).Value = 37500  # J59
$ws.Cells.Item(59, 12).Value = 37500  # L59
$ws.Cells.Item(59, 14).Value = -39194  # N59
$ws.Cells.Item(107, 8).Value = 5091.6  # H107
$ws.Cells.Item(107, 9).Value = 737.3182  # I107
$ws.Cells.Item(107, 11).Value = 737.3182  # K107
$ws.Cells.Item(107, 13).Value = 1182.6818  # M107

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2632.6  # H16
$ws.Cells.Item(16, 9).Value = 700  # I16
$ws.Cells.Item(16, 11).Value = 700  # K16
$ws.Cells.Item(16, 13).Value = -413  # M16
$ws.Cells.Item(106, 8).Value = 0  # H106
$ws.Cells.Item(106, 10).Value = 0  # J106
$ws.Cells.Item(106, 12).Value = ""  # clear L106
$ws.Cells.Item(106, 14).Value = 0  # N106
$ws.Cells.Item(113, 8).Value = 2632.6  # H113
$ws.Cells.Item(113, 9).Value = 700  # I113
$ws.Cells.Item(113, 11).Value = 700  # K113
$ws.Cells.Item(113, 13).Value = 1470  # M113
$ws.Cells.Item(134, 8).Value = 1204.7115  # H134
$ws.Cells.Item(134, 9).Value = 1049.9736  # I134
$ws.Cells.Item(134, 10).Value = 1624.7142  # J134
$ws.Cells.Item(134, 11).Value = 3149.9208  # K134
$ws.Cells.Item(134, 12).Value = 4874.142599999999  # L134
$ws.Cells.Item(134, 13).Value = -614.9207999999999  # M134
$ws.Cells.Item(134, 14).Value = -9944.142599999999  # N134

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(88, 8).Value = 7928.5713  # H88
$ws.Cells.Item(88, 9).Value = 0  # I88
$ws.Cells.Item(88, 10).Value = 7928.5713  # J88
$ws.Cells.Item(88, 11).Value = 0  # K88
$ws.Cells.Item(88, 12).Value = ""  # clear L88
$ws.Cells.Item(88, 13).Value = 23785.7139  # M88
$ws.Cells.Item(88, 14).Value = -24641.7139  # N88
$ws.Cells.Item(91, 8).Value = 7928.5713  # H91
$ws.Cells.Item(91, 9).Value = 0  # I91
$ws.Cells.Item(91, 10).Value = 7928.5713  # J91
$ws.Cells.Item(91, 11).Value = 0  # K91
$ws.Cells.Item(91, 12).Value = ""  # clear L91
$ws.Cells.Item(91, 13).Value = 23785.7139  # M91
$ws.Cells.Item(91, 14).Value = -26749.7139  # N91
$ws.Cells.Item(113, 8).Value = 568  # H113
$ws.Cells.Item(113, 9).Value = 668.5714  # I113
$ws.Cells.Item(113, 10).Value = 521.06665  # J113
$ws.Cells.Item(113, 11).Value = 2005.7142  # K113
$ws.Cells.Item(113, 12).Value = 1563.19995  # L113
$ws.Cells.Item(113, 13).Value = 164.2857999999999  # M113
$ws.Cells.Item(113, 14).Value = -5903.19995  # N113
$ws.Cells.Item(129, 8).Value = 55290.26  # H129
$ws.Cells.Item(129, 9).Value = 1220  # I129
$ws.Cells.Item(129, 10).Value = 94614.09  # J129
$ws.Cells.Item(129, 11).Value = 3660  # K129
$ws.Cells.Item(129, 12).Value = 283842.27  # L129
$ws.Cells.Item(129, 13).Value = 1340  # M129
$ws.Cells.Item(129, 14).Value = -293842.27  # N129

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(34, 8).Value = 0  # H34
$ws.Cells.Item(34, 9).Value = 0  # I34
$ws.Cells.Item(34, 11).Value = 0  # K34
$ws.Cells.Item(34, 13).Value = ""  # clear M34
$ws.Cells.Item(70, 8).Value = 4111.9165  # H70
$ws.Cells.Item(70, 9).Value = 3922.6296  # I70
$ws.Cells.Item(70, 11).Value = 3922.6296  # K70
$ws.Cells.Item(70, 13).Value = -3652.6296  # M70
$ws.Cells.Item(73, 8).Value = 4111.9165  # H73
$ws.Cells.Item(73, 9).Value = 3922.6296  # I73
$ws.Cells.Item(73, 11).Value = 3922.6296  # K73
$ws.Cells.Item(73, 13).Value = -2986.6296  # M73
$ws.Cells.Item(74, 8).Value = 59543.668  # H74
$ws.Cells.Item(74, 10).Value = 59543.668  # J74
$ws.Cells.Item(74, 12).Value = 59543.668  # L74
$ws.Cells.Item(74, 14).Value = -61415.668  # N74
$ws.Cells.Item(75, 8).Value = 8000  # H75
$ws.Cells.Item(75, 10).Value = 8000  # J75
$ws.Cells.Item(75, 12).Value = 8000  # L75
$ws.Cells.Item(75, 14).Value = -9748  # N75
$ws.Cells.Item(76, 8).Value = 0  # H76
$ws.Cells.Item(76, 9).Value = 0  # I76
$ws.Cells.Item(76, 11).Value = 0  # K76
$ws.Cells.Item(76, 13).Value = ""  # clear M76
$ws.Cells.Item(77, 8).Value = 59543.668  # H77
$ws.Cells.Item(77, 10).Value = 59543.668  # J77
$ws.Cells.Item(77, 12).Value = 178631.004  # L77
$ws.Cells.Item(77, 14).Value = -187991.004  # N77
$ws.Cells.Item(78, 8).Value = 8000  # H78
$ws.Cells.Item(78, 10).Value = 8000  # J78
$ws.Cells.Item(78, 12).Value = 24000  # L78
$ws.Cells.Item(78, 14).Value = -32736  # N78
$ws.Cells.Item(79, 8).Value = 0  # H79
$ws.Cells.Item(79, 9).Value = 0  # I79
$ws.Cells.Item(79, 11).Value = 0  # K79
$ws.Cells.Item(79, 13).Value = ""  # clear M79
$ws.Cells.Item(97, 8).Value = 1244.0435  # H97
$ws.Cells.Item(97, 9).Value = 1122.25  # I97
$ws.Cells.Item(97, 11).Value = 1122.25  # K97
$ws.Cells.Item(97, 13).Value = -626.25  # M97

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3369003.2  # H40
$ws.Cells.Item(40, 9).Value = 3369003.2  # I40
$ws.Cells.Item(40, 11).Value = 3369003.2  # K40
$ws.Cells.Item(40, 13).Value = -3368867.2  # M40
$ws.Cells.Item(82, 8).Value = 1981.375  # H82
$ws.Cells.Item(82, 9).Value = 2160.2  # I82
$ws.Cells.Item(82, 11).Value = 2160.2  # K82
$ws.Cells.Item(82, 13).Value = -1799.2  # M82
$ws.Cells.Item(85, 8).Value = 1981.375  # H85
$ws.Cells.Item(85, 9).Value = 2160.2  # I85
$ws.Cells.Item(85, 11).Value = 2160.2  # K85
$ws.Cells.Item(85, 13).Value = -912.1999999999998  # M85
$ws.Cells.Item(122, 8).Value = 1700  # H122
$ws.Cells.Item(122, 9).Value = 1600  # I122
$ws.Cells.Item(122, 10).Value = 2200  # J122
$ws.Cells.Item(122, 11).Value = 4800  # K122
$ws.Cells.Item(122, 12).Value = 6600  # L122
$ws.Cells.Item(122, 13).Value = -2350  # M122
$ws.Cells.Item(122, 14).Value = -11500  # N122
$ws.Cells.Item(136, 8).Value = 3235.743  # H136
$ws.Cells.Item(136, 9).Value = 1159.1212  # I136
$ws.Cells.Item(136, 10).Value = 37500  # J136
$ws.Cells.Item(136, 11).Value = 3477.3636  # K136
$ws.Cells.Item(136, 12).Value = 112500  # L136
$ws.Cells.Item(136, 13).Value = -927.3636000000001  # M136
$ws.Cells.Item(136, 14).Value = -117600  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 0  # H92
$ws.Cells.Item(92, 10).Value = 0  # J92
$ws.Cells.Item(92, 12).Value = ""  # clear L92
$ws.Cells.Item(92, 14).Value = 0  # N92
$ws.Cells.Item(122, 8).Value = 917.3333  # H122
$ws.Cells.Item(122, 9).Value = 800.8  # I122
$ws.Cells.Item(122, 10).Value = 1500  # J122
$ws.Cells.Item(122, 11).Value = 2402.4  # K122
$ws.Cells.Item(122, 12).Value = 4500  # L122
$ws.Cells.Item(122, 13).Value = 47.60000000000036  # M122
$ws.Cells.Item(122, 14).Value = -9400  # N122
$ws.Cells.Item(133, 8).Value = 76223  # H133
$ws.Cells.Item(133, 10).Value = 76223  # J133
$ws.Cells.Item(133, 12).Value = 76223  # L133
$ws.Cells.Item(133, 14).Value = -86343  # N133
$ws.Cells.Item(136, 8).Value = 599.5357  # H136
$ws.Cells.Item(136, 9).Value = 359.94736  # I136
$ws.Cells.Item(136, 10).Value = 1105.3334  # J136
$ws.Cells.Item(136, 11).Value = 1079.84208  # K136
$ws.Cells.Item(136, 12).Value = 3316.0002  # L136
$ws.Cells.Item(136, 13).Value = 1470.15792  # M136
$ws.Cells.Item(136, 14).Value = -8416.0002  # N136
